# se modifica data para hacer smoke en preProd
$wb = $excel.ActiveWorkbook

# --- DatosCuenta sheet (SmokPreProdOcho -> SmokPreProdNueve data set) ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokPreProdNueve"
$wsCuenta.Range("B2").Value = "SmokeNamePreProdNueve"
$wsCuenta.Range("C2").Value = 27100117
$wsCuenta.Range("D2").Value = 119

# --- DatosHogar sheet ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 638

# --- DatosMotor sheet (SMP019 -> SMP020 plate/engine/chassis data) ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP020"
$wsMotor.Range("B2").Value = "ABC12SSMP020"
$wsMotor.Range("C2").Value = "ZAZ123SSMP020"
$wsMotor.Range("A2:C2").Select() | Out-Null

# --- DatosAP sheet ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200120
$wsAP.Range("E11").Select() | Out-Null

# Make DatosCuenta the active (selected) sheet/tab
$wsCuenta.Activate() | Out-Null
$wsCuenta.Range("D2").Select() | Out-Null
